# Financials update: add the new FY2018 (period ending 2018-12-31) column
# in front of the existing year columns on the "MAN" sheet, across the
# Income Statement, Balance Sheet and Cash Flow Statement blocks, and
# restate a handful of FY2017/FY2016/FY2015 figures that changed together
# with this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D - this shifts the existing D:K data (the
#    2017..2011 year columns) one column to the right (E:L), carrying
#    over cell styles/number formats automatically.
$ws.Range("D1").EntireColumn.Insert()

# Keep the new column's width consistent with its neighbours (the other
# year columns), since Excel leaves it at a default width after insert.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 2) Populate the new D column (2018-12-31 figures) and patch the small
#    number of restated figures (2017 / 2016 columns) called out in this
#    update. Column letters below refer to the POST-insert layout.

# --- Income Statement (header row 7, data rows 8-35) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 21991200
$ws.Range("D9").Value = 18412200
$ws.Range("D10").Value = 3579000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 39300
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 21194500
$ws.Range("E17").Value = 20245100
$ws.Range("F17").Value = 18908600
$ws.Range("D18").Value = 796700
$ws.Range("E18").Value = 789200
$ws.Range("F18").Value = 745500
$ws.Range("D20").Value = 5000
$ws.Range("E20").Value = -2500
$ws.Range("F20").Value = 5300
$ws.Range("D21").Value = 887500
$ws.Range("E21").Value = 871100
$ws.Range("F21").Value = 836100
$ws.Range("D22").Value = 47000
$ws.Range("E22").Value = 49400
$ws.Range("F22").Value = 49500
$ws.Range("D23").Value = 754700
$ws.Range("D24").Value = 194800
$ws.Range("E24").Value = -255300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 559900
$ws.Range("E26").Value = 992600
$ws.Range("D27").Value = 559900
$ws.Range("E27").Value = 992600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -3200
$ws.Range("E29").Value = -447200
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -5000
$ws.Range("E32").Value = 2500
$ws.Range("F32").Value = -5300
$ws.Range("D33").Value = 556700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 556700

# --- Balance Sheet (header row 38, data rows 39-77) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 591900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 5276100
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 129100
$ws.Range("D46").Value = 5997100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 152600
$ws.Range("D49").Value = 1543400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 826700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 8519800
$ws.Range("D57").Value = 2266700
$ws.Range("D58").Value = 50100
$ws.Range("D59").Value = 1859100
$ws.Range("D60").Value = 4175900
$ws.Range("D61").Value = 1025300
$ws.Range("D62").Value = 620100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 5894900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3157700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2624900
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (header row 80, data rows 81-102) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 556700
$ws.Range("D83").Value = 85800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 483100
$ws.Range("D91").Value = -64700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -54900
$ws.Range("D96").Value = -127300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -484900
$ws.Range("D101").Value = -40400
$ws.Range("D102").Value = -97100
